$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "60.922.52"
$ws.Range("E2").Value = "  +2.76%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.608.67"
$ws.Range("E3").Value = "  +1.11%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
Set-TextValue $ws.Range("D5") "579.10"
$ws.Range("E5").Value = "  +4.23%  "

# Row 6
Set-TextValue $ws.Range("D6") "143.88"
$ws.Range("E6").Value = "  +1.20%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.996"
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("E8").Value = "  +0.46%  "

# Row 9
Set-TextValue $ws.Range("D9") "2.633.36"
$ws.Range("E9").Value = "  +1.91%  "

# Row 10
Set-TextValue $ws.Range("D10") "6.54"
$ws.Range("E10").Value = "  -2.82%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.107"
$ws.Range("E11").Value = "  +2.40%  "

# Row 12
$ws.Range("E12").Value = "  -5.44%  "

# Row 13
$ws.Range("E13").Value = "  +5.49%  "

# Row 14
Set-TextValue $ws.Range("D14") "3.072.21"
$ws.Range("E14").Value = "  +1.27%  "

# Row 15
Set-TextValue $ws.Range("D15") "60.877.67"
$ws.Range("E15").Value = "  +2.77%  "

# Row 16
Set-TextValue $ws.Range("D16") "23.47"
$ws.Range("E16").Value = "  +1.89%  "

# Row 17
$ws.Range("E17").Value = "  +4.57%  "

# Row 18
Set-TextValue $ws.Range("D18") "2.620.76"
$ws.Range("E18").Value = "  +1.40%  "

# Row 19
Set-TextValue $ws.Range("D19") "11.31"
$ws.Range("E19").Value = "  +9.27%  "

# Row 20
Set-TextValue $ws.Range("D20") "4.68"
$ws.Range("E20").Value = "  +2.91%  "

# Row 21
Set-TextValue $ws.Range("D21") "350.64"
$ws.Range("E21").Value = "  +3.89%  "

# Row 22
Set-TextValue $ws.Range("D22") "6.97"
$ws.Range("E22").Value = "  +7.92%  "

# Row 23
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
Set-TextValue $ws.Range("D24") "0.520"
$ws.Range("E24").Value = "  +9.34%  "

# Row 25
Set-TextValue $ws.Range("D25") "63.35"
$ws.Range("E25").Value = "  +1.36%  "

# Row 26
$ws.Range("E26").Value = "  -0.55%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.162"
$ws.Range("E27").Value = "  +1.61%  "

# Row 28
Set-TextValue $ws.Range("D28") "7.99"
$ws.Range("E28").Value = "  +7.77%  "

# Row 29
Set-TextValue $ws.Range("D29") "0.0₃0806"
$ws.Range("E29").Value = "  +3.96%  "

# Row 30
$ws.Range("E30").Value = "  +9.74%  "

# Row 31
$ws.Range("E31").Value = "  +2.50%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.997"
$ws.Range("E32").Value = "  -0.11%  "

# Row 33
Set-TextValue $ws.Range("D33") "162.49"
$ws.Range("E33").Value = "  +2.32%  "

# Row 34
Set-TextValue $ws.Range("D34") "19.59"
$ws.Range("E34").Value = "  +2.62%  "

# Row 35
$ws.Range("E35").Value = "  +15.65%  "

# Row 36
$ws.Range("E36").Value = "  +5.41%  "

# Row 37
Set-TextValue $ws.Range("D37") "1.25"
$ws.Range("E37").Value = "  +6.50%  "

# Row 38
Set-TextValue $ws.Range("D38") "1.63"
$ws.Range("E38").Value = "  +9.56%  "

# Row 39
Set-TextValue $ws.Range("D39") "37.96"
$ws.Range("E39").Value = "  +1.61%  "

# Row 40
$ws.Range("E40").Value = "  +5.79%  "

# Row 41
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D41") "306.73"
$ws.Range("E41").Value = "  +5.80%  "

# Row 42
$ws.Range("B42").Value = "SuiNetwork"
$ws.Range("C42").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D42") "0.854"
$ws.Range("E42").Value = "  +0.02%  "

# Row 43
Set-TextValue $ws.Range("D43") "134.58"
$ws.Range("E43").Value = "  -1.83%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D44") "20.59"
$ws.Range("E44").Value = "  +10.14%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D45") "20.04"
$ws.Range("E45").Value = "  +5.87%  "

# Row 46
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D46") "5.06"
$ws.Range("E46").Value = "  +11.96%  "

# Row 47
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D47") "0.995"
$ws.Range("E47").Value = "  -0.37%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D48") "0.0987"
$ws.Range("E48").Value = "  +1.32%  "

# Row 49
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D49") "0.608"
$ws.Range("E49").Value = "  +2.58%  "

# Row 50
$ws.Range("B50").Value = "Hedera"
$ws.Range("C50").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D50") "0.0553"
$ws.Range("E50").Value = "  +4.28%  "

# Row 51
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D51") "0.0243"
$ws.Range("E51").Value = "  +3.90%  "
